$d = $word.ActiveDocument

function Set-BoldForText($doc, [string]$text) {
    $range = $doc.Content.Duplicate
    $found = $range.Find.Execute($text, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $range.Font.Bold = 1
    }
}

# 1) Bold the name "Holly Dickson"
Set-BoldForText $d "Holly Dickson"

# 2) Bold "Experiência de trabalho" heading run
Set-BoldForText $d "Experiência de trabalho"

# 3) Bold the job title run "Designer júnior de animação" (before the text itself changes)
Set-BoldForText $d "Designer júnior de animação"

# 4) Change job title text "Designer júnior de animação" -> "Designer de animação junior"
$d.Content.Find.Execute("Designer júnior de animação", $true, $true, $false, $false, $false, $true, 1, $false, "Designer de animação junior", 2) | Out-Null

# 5) Bold "Estagiário de animação" heading run
Set-BoldForText $d "Estagiário de animação"

# 6) Change text "Estagiário de animação" -> "Estagiário de Animação" (capitalize A)
$d.Content.Find.Execute("Estagiário de animação", $true, $true, $false, $false, $false, $true, 1, $false, "Estagiário de Animação", 2) | Out-Null

# 7) Bold the degree heading run "Bacharelado em Belas Artes em Animação"
Set-BoldForText $d "Bacharelado em Belas Artes em Animação"
